$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 13:45:48"
$ws.Cells.Item(3,1).Value = "Total filas: 203"
$ws.Cells.Item(22,3).Value = "14_ABASTO"
$ws.Cells.Item(23,3).Value = "215C_EL PATO"
$ws.Cells.Item(68,3).Value = "215A_EL PATO"
$ws.Cells.Item(69,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(70,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(76,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(77,3).Value = "16_SANTA ANA"
$ws.Cells.Item(82,3).Value = "17_ROMERO"
$ws.Cells.Item(84,3).Value = "16_SANTA ANA"
$ws.Cells.Item(98,1).Value = "09:01:18"
$ws.Cells.Item(98,3).Value = "10_OLMOS"
$ws.Cells.Item(98,4).Value = 70
$ws.Cells.Item(99,1).Value = "08:41:16"
$ws.Cells.Item(99,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(99,4).Value = 90
$ws.Cells.Item(116,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(117,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(141,1).Value = "11:07:42"
$ws.Cells.Item(141,3).Value = "14_ABASTO"
$ws.Cells.Item(141,4).Value = 59
$ws.Cells.Item(142,1).Value = "11:54:47"
$ws.Cells.Item(142,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(142,4).Value = 12
$ws.Cells.Item(143,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(163,3).Value = "14_ABASTO"
$ws.Cells.Item(164,3).Value = "15_ABASTO"
$ws.Cells.Item(179,1).Value = "13:45:48"
$ws.Cells.Item(179,3).Value = "17_ROMERO"
$ws.Cells.Item(179,4).Value = 1
$ws.Cells.Item(180,1).Value = "13:45:48"
$ws.Cells.Item(180,3).Value = "16_SANTA ANA"
$ws.Cells.Item(180,4).Value = 1
$ws.Cells.Item(181,1).Value = "13:45:48"
$ws.Cells.Item(181,3).Value = "215A_EL PATO"
$ws.Cells.Item(181,4).Value = 5
$ws.Cells.Item(182,1).Value = "13:45:48"
$ws.Cells.Item(182,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(182,4).Value = 5
$ws.Cells.Item(183,1).Value = "13:45:48"
$ws.Cells.Item(183,2).Value = "13:55"
$ws.Cells.Item(183,3).Value = "225_GOMEZ"
$ws.Cells.Item(183,4).Value = 10
$ws.Cells.Item(184,1).Value = "13:45:48"
$ws.Cells.Item(184,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(184,4).Value = 11
$ws.Cells.Item(185,2).Value = "13:56"
$ws.Cells.Item(185,3).Value = "225_GOMEZ"
$ws.Cells.Item(185,4).Value = 57
$ws.Cells.Item(186,1).Value = "13:45:48"
$ws.Cells.Item(186,2).Value = "14:04"
$ws.Cells.Item(186,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(186,4).Value = 19
$ws.Cells.Item(187,1).Value = "13:45:48"
$ws.Cells.Item(187,2).Value = "14:04"
$ws.Cells.Item(187,3).Value = "17_ROMERO"
$ws.Cells.Item(187,4).Value = 19
$ws.Cells.Item(188,1).Value = "13:45:48"
$ws.Cells.Item(188,2).Value = "14:06"
$ws.Cells.Item(188,3).Value = "16_SANTA ANA"
$ws.Cells.Item(188,4).Value = 21
$ws.Cells.Item(189,1).Value = "13:45:48"
$ws.Cells.Item(189,2).Value = "14:12"
$ws.Cells.Item(189,3).Value = "15_ABASTO"
$ws.Cells.Item(189,4).Value = 27
$ws.Cells.Item(190,1).Value = "13:45:48"
$ws.Cells.Item(190,2).Value = "14:16"
$ws.Cells.Item(190,3).Value = "27_EL RETIRO"
$ws.Cells.Item(190,4).Value = 31
$ws.Cells.Item(191,1).Value = "13:45:48"
$ws.Cells.Item(191,2).Value = "14:17"
$ws.Cells.Item(191,3).Value = "27_EL RETIRO"
$ws.Cells.Item(191,4).Value = 78
$ws.Cells.Item(192,1).Value = "13:45:48"
$ws.Cells.Item(192,2).Value = "14:19"
$ws.Cells.Item(192,3).Value = "215C_EL PATO"
$ws.Cells.Item(192,4).Value = 34
$ws.Cells.Item(193,2).Value = "14:20"
$ws.Cells.Item(193,3).Value = "215C_EL PATO"
$ws.Cells.Item(193,4).Value = 81
$ws.Cells.Item(194,1).Value = "13:45:48"
$ws.Cells.Item(194,2).Value = "14:21"
$ws.Cells.Item(194,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(194,4).Value = 36
$ws.Cells.Item(194,5).Value = "LP1912"
$ws.Cells.Item(195,1).Value = "12:59:25"
$ws.Cells.Item(195,2).Value = "14:42"
$ws.Cells.Item(195,3).Value = "14_ABASTO"
$ws.Cells.Item(195,4).Value = 103
$ws.Cells.Item(195,5).Value = "LP1912"
$ws.Cells.Item(196,1).Value = "13:45:48"
$ws.Cells.Item(196,2).Value = "14:44"
$ws.Cells.Item(196,3).Value = "14_ABASTO"
$ws.Cells.Item(196,4).Value = 59
$ws.Cells.Item(196,5).Value = "LP1912"
$ws.Cells.Item(197,1).Value = "13:45:48"
$ws.Cells.Item(197,2).Value = "14:56"
$ws.Cells.Item(197,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(197,4).Value = 71
$ws.Cells.Item(197,5).Value = "LP1912"
$ws.Cells.Item(198,1).Value = "13:45:48"
$ws.Cells.Item(198,2).Value = "14:58"
$ws.Cells.Item(198,3).Value = "215B_EL PATO"
$ws.Cells.Item(198,4).Value = 73
$ws.Cells.Item(198,5).Value = "LP1912"
$ws.Cells.Item(199,1).Value = "13:45:48"
$ws.Cells.Item(199,2).Value = "15:00"
$ws.Cells.Item(199,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(199,4).Value = 75
$ws.Cells.Item(199,5).Value = "LP1912"
$ws.Cells.Item(200,1).Value = "13:45:48"
$ws.Cells.Item(200,2).Value = "15:05"
$ws.Cells.Item(200,3).Value = "10_OLMOS"
$ws.Cells.Item(200,4).Value = 80
$ws.Cells.Item(200,5).Value = "LP1912"
$ws.Cells.Item(201,1).Value = "13:45:48"
$ws.Cells.Item(201,2).Value = "15:13"
$ws.Cells.Item(201,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(201,4).Value = 88
$ws.Cells.Item(201,5).Value = "LP1912"
$ws.Cells.Item(202,1).Value = "13:45:48"
$ws.Cells.Item(202,2).Value = "15:20"
$ws.Cells.Item(202,3).Value = "15_ABASTO"
$ws.Cells.Item(202,4).Value = 95
$ws.Cells.Item(202,5).Value = "LP1912"
$ws.Cells.Item(203,1).Value = "13:45:48"
$ws.Cells.Item(203,2).Value = "15:22"
$ws.Cells.Item(203,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(203,4).Value = 97
$ws.Cells.Item(203,5).Value = "LP1912"
$ws.Cells.Item(204,1).Value = "13:45:48"
$ws.Cells.Item(204,2).Value = "15:32"
$ws.Cells.Item(204,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(204,4).Value = 107
$ws.Cells.Item(204,5).Value = "LP1912"
$ws.Cells.Item(205,1).Value = "13:45:48"
$ws.Cells.Item(205,2).Value = "15:34"
$ws.Cells.Item(205,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(205,4).Value = 109
$ws.Cells.Item(205,5).Value = "LP1912"
$ws.Cells.Item(206,1).Value = "13:45:48"
$ws.Cells.Item(206,2).Value = "15:37"
$ws.Cells.Item(206,3).Value = "10_OLMOS"
$ws.Cells.Item(206,4).Value = 112
$ws.Cells.Item(206,5).Value = "LP1912"
$ws.Cells.Item(207,1).Value = "13:45:48"
$ws.Cells.Item(207,2).Value = "15:38"
$ws.Cells.Item(207,3).Value = "215A_EL PATO"
$ws.Cells.Item(207,4).Value = 113
$ws.Cells.Item(207,5).Value = "LP1912"
$ws.Cells.Item(208,1).Value = "13:45:48"
$ws.Cells.Item(208,2).Value = "15:42"
$ws.Cells.Item(208,3).Value = "14_ABASTO"
$ws.Cells.Item(208,4).Value = 117
$ws.Cells.Item(208,5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 13:45:48"
$ws.Cells.Item(3,1).Value = "Total filas: 23"
$ws.Cells.Item(24,1).Value = "13:45:48"
$ws.Cells.Item(24,4).Value = 5
$ws.Cells.Item(25,1).Value = "13:45:48"
$ws.Cells.Item(25,4).Value = 34
$ws.Cells.Item(27,1).Value = "13:45:48"
$ws.Cells.Item(27,2).Value = "14:58"
$ws.Cells.Item(27,3).Value = "215B_EL PATO"
$ws.Cells.Item(27,4).Value = 73
$ws.Cells.Item(27,5).Value = "LP1912"
$ws.Cells.Item(28,1).Value = "13:45:48"
$ws.Cells.Item(28,2).Value = "15:38"
$ws.Cells.Item(28,3).Value = "215A_EL PATO"
$ws.Cells.Item(28,4).Value = 113
$ws.Cells.Item(28,5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 13:45:48"
$ws.Cells.Item(3,1).Value = "Total filas: 30"
$ws.Cells.Item(32,1).Value = "13:45:48"
$ws.Cells.Item(32,4).Value = 24
$ws.Cells.Item(33,1).Value = "13:45:48"
$ws.Cells.Item(33,2).Value = "14:52"
$ws.Cells.Item(33,4).Value = 67
$ws.Cells.Item(34,1).Value = "12:59:25"
$ws.Cells.Item(34,2).Value = "14:53"
$ws.Cells.Item(34,3).Value = "215D_LA PLATA"
$ws.Cells.Item(34,4).Value = 114
$ws.Cells.Item(34,5).Value = "L6203"
$ws.Cells.Item(35,1).Value = "13:45:48"
$ws.Cells.Item(35,2).Value = "15:34"
$ws.Cells.Item(35,3).Value = "215A_LA PLATA"
$ws.Cells.Item(35,4).Value = 109
$ws.Cells.Item(35,5).Value = "L6173"
